# Refresh the cryptocurrency price/volume snapshot (and fix the rank-order
# swaps for WEMIXToken/TrustWalletToken and InjectiveProtocol/Aave).
#
# Column D ("Price") cells are plain text in the source workbook (e.g.
# "35.363.67", "0.730", "5.00"). A bare assignment like
# $ws.Range('D16').Value = '5.00' would let Excel auto-coerce the text into
# the number 5, which both changes the cell's stored type and silently
# drops the trailing zero. Prefixing the literal with an apostrophe forces
# Excel to keep it as text, matching the original formatting exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''35.398.16'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '''1.911.76'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '''0.731'
$ws.Range('E5').Value = '  +10.81%  '
$ws.Range('D6').Value = '''256.32'
$ws.Range('E6').Value = '  +4.30%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').Value = '''41.08'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '''0.371'
$ws.Range('E9').Value = '  +7.19%  '
$ws.Range('D10').Value = '''53.19'
$ws.Range('E10').Value = '  +0.96%  '
$ws.Range('D11').Value = '''0.0761'
$ws.Range('E11').Value = '  +6.42%  '
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '''2.186.63'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = '''12.96'
$ws.Range('E14').Value = '  +5.74%  '
$ws.Range('D15').Value = '''0.735'
$ws.Range('E15').Value = '  +6.01%  '
$ws.Range('D16').Value = '''5.00'
$ws.Range('E16').Value = '  +4.44%  '
$ws.Range('D17').Value = '''1.913.32'
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('D18').Value = '''35.369.26'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').Value = '''75.38'
$ws.Range('E19').Value = '  +4.65%  '
$ws.Range('D20').Value = '''0.0₃0848'
$ws.Range('E20').Value = '  +4.07%  '
$ws.Range('D21').Value = '''245.47'
$ws.Range('E21').Value = '  +2.48%  '
$ws.Range('D22').Value = '''13.13'
$ws.Range('E22').Value = '  +6.10%  '
$ws.Range('D23').Value = '''5.16'
$ws.Range('E23').Value = '  +7.49%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  +6.88%  '
$ws.Range('E26').Value = '  +3.23%  '
$ws.Range('D27').Value = '''166.32'
$ws.Range('E27').Value = '  -2.32%  '
$ws.Range('D28').Value = '''8.77'
$ws.Range('E28').Value = '  +4.18%  '
$ws.Range('D29').Value = '''18.85'
$ws.Range('E29').Value = '  +3.12%  '
$ws.Range('E30').Value = '  +5.10%  '
$ws.Range('D31').Value = '''4.128.99'
$ws.Range('E31').Value = '  -1.35%  '
$ws.Range('E32').Value = '  +6.48%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '''2.00'
$ws.Range('E33').Value = '  +14.34%  '
$ws.Range('B34').Value = 'TrustWalletToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D34').Value = '''1.65'
$ws.Range('E34').Value = '  +24.71%  '
$ws.Range('D35').Value = '''0.0592'
$ws.Range('E35').Value = '  +5.83%  '
$ws.Range('E36').Value = '  +5.42%  '
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('D38').Value = '''0.916'
$ws.Range('E38').Value = '  -2.56%  '
$ws.Range('D39').Value = '''2.05'
$ws.Range('E39').Value = '  +1.62%  '
$ws.Range('E40').Value = '  +6.57%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = '''17.14'
$ws.Range('E41').Value = '  +7.11%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '''97.59'
$ws.Range('E42').Value = '  +9.33%  '
$ws.Range('E43').Value = '  +3.41%  '
$ws.Range('E44').Value = '  +1.31%  '
$ws.Range('E45').Value = '  +5.19%  '
$ws.Range('D46').Value = '''1.343.27'
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('E47').Value = '  +0.83%  '
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').Value = '''45.34'
$ws.Range('E50').Value = '  -8.54%  '
$ws.Range('D51').Value = '''0.0756'
$ws.Range('E51').Value = '  +6.80%  '
